$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header / label cells (order matters for shared-string table order)
$ws.Range("E7").Value = "logx"
$ws.Range("F7").Value = "logy"
$ws.Range("P8").Value = "rapport"
$ws.Range("X8").Value = "constatnes"

# CSV-formatted (period decimal separator) literal text values, mirroring
# column T's CONCAT results but pasted as plain text values in column X
$ws.Range("X1").Value = "4096.0 0.00564897060394285"
$ws.Range("X2").Value = "32768.0 0.0102252244949341"
$ws.Range("X3").Value = "262144.0 0.0960252761840819"
$ws.Range("X4").Value = "2097152.0 0.768564534187317"
$ws.Range("X5").Value = "16777216.0 6.3871531009674"
$ws.Range("X6").Value = "134217728.0 87.6655602693557"

# Update the saved view state (scroll position + active selection)
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("U11").Select()
